$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Username was re-capitalised ("Hpetrov" -> "HPetrov") as part of the
# password-changing feature / self-evaluation protocol update.
$ws.Range("C4").Value = "HPetrov"

# GitHub activity numbers filled in.
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 34

# Basic-options scores filled in (Change Password / Edit User Profile).
$ws.Range("C29").Value = 5
$ws.Range("C30").Value = 5

# Move the active selection/cursor to F8 (was G26, with the view
# previously scrolled so A21 was the top-left cell).
$ws.Range("F8").Select()
